$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "The candidate has expertise in JavaScript, ReactJS, and Python, which are common required skills. However, the candidate lacks experience in CSS and NextJS, which are required for the job. The candidate's projects demonstrate strong technical skills, especially in AI, but lack specific frontend development experience."
$ws.Range("E2").Value = "The candidate shows a high interest in AI and demonstrates good teamwork and problem-solving skills. The willingness to learn a new language and adapt to a new culture shows adaptability. However, the lack of experience in working alone and the preference for a stable job after graduation might impact long-term commitment."

$ws.Range("C3").Value = "The candidate has experience with required skills ReactJS, JavaScript, and CSS, along with additional skills like NodeJS, ExpressJS, and Firebase. The candidate lacks experience with NextJS, but the overall experience and project work demonstrate a good fit for the role."
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = "The applicant demonstrates strong communication, teamwork, adaptability, and problem-solving skills. Their willingness to adapt to Japanese work culture and learn the language shows a positive attitude towards working in Japan. Overall, the candidate's responses align with the company's soft skill requirements and willingness to relocate, warranting a high personality score."

$ws.Range("C4").Value = "The candidate lacks the required skills of Frontend Development, CSS, and NextJS which are crucial for the job role. The candidate also lacks experience in collaborating with UX and graphic designers, as mentioned in the job description."
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "The applicant shows strong interest in AI and Japanese work culture and is open to adapting to a new environment. He also demonstrates good teamwork and problem-solving skills, which are essential for the internship."

$wb.Save()
